# Fixed extraction and algorithm to generate uncapped
# Updates column E ("imgVW") values on both the "listing" and "detail"
# worksheets for rows 2-17, replacing the capped values with the
# uncapped ones produced by the corrected algorithm.

$wb = $excel.ActiveWorkbook

# Row -> new value mapping for each sheet's column E (imgVW)
$listingValues = @{
    2  = 45
    3  = 45
    4  = 45
    5  = 45
    6  = 45
    7  = 45
    8  = 45
    9  = 23
    10 = 45
    11 = 23
    12 = 23
    13 = 45
    14 = 45
    15 = 45
    16 = 23
    17 = 45
}

$detailValues = @{
    2  = 45
    3  = 45
    4  = 45
    5  = 45
    6  = 45
    7  = 45
    8  = 45
    9  = 40
    10 = 45
    11 = 38
    12 = 38
    13 = 45
    14 = 45
    15 = 45
    16 = 38
    17 = 45
}

$wsListing = $wb.Worksheets.Item("listing")
foreach ($row in $listingValues.Keys) {
    $wsListing.Cells.Item($row, 5).Value = $listingValues[$row]
}

$wsDetail = $wb.Worksheets.Item("detail")
foreach ($row in $detailValues.Keys) {
    $wsDetail.Cells.Item($row, 5).Value = $detailValues[$row]
}
